$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (the "Förändrad" date column) for rows 2 through 45
# from serial date 45743 (2025-03-27) to 45744 (2025-03-28).
$ws.Range("C2:C45").Value = 45744
